{"js": "// Replace the 15 lattice-multiplication exercise cells (5 rows x 3 cols)\n// with the new set of problems, keeping the table/cell/run formatting intact.\n// Each cell holds a single run with 5 \"lines\" separated by <w:br/> (which the\n// Word JS text model exposes as \\u000b \"vertical tab\").\nconst newCellText = [\n  \"60 x 73\\u000b  7    3\\u000b  ----\\u000b6|    |\\u000b0|    |\",\n  \"21 x 75\\u000b  7    5\\u000b  ----\\u000b2|    |\\u000b1|    |\",\n  \"58 x 47\\u000b  4    7\\u000b  ----\\u000b5|    |\\u000b8|    |\",\n  \"60 x 86\\u000b  8    6\\u000b  ----\\u000b6|    |\\u000b0|    |\",\n  \"56 x 42\\u000b  4    2\\u000b  ----\\u000b5|    |\\u000b6|    |\",\n  \"29 x 36\\u000b  3    6\\u000b  ----\\u000b2|    |\\u000b9|    |\",\n  \"20 x 89\\u000b  8    9\\u000b  ----\\u000b2|    |\\u000b0|    |\",\n  \"94 x 85\\u000b  8    5\\u000b  ----\\u000b9|    |\\u000b4|    |\",\n  \"88 x 40\\u000b  4    0\\u000b  ----\\u000b8|    |\\u000b8|    |\",\n  \"86 x 34\\u000b  3    4\\u000b  ----\\u000b8|    |\\u000b6|    |\",\n  \"15 x 55\\u000b  5    5\\u000b  ----\\u000b1|    |\\u000b5|    |\",\n  \"96 x 98\\u000b  9    8\\u000b  ----\\u000b9|    |\\u000b6|    |\",\n  \"43 x 24\\u000b  2    4\\u000b  ----\\u000b4|    |\\u000b3|    |\",\n  \"15 x 78\\u000b  7    8\\u000b  ----\\u000b1|    |\\u000b5|    |\",\n  \"49 x 42\\u000b  4    2\\u000b  ----\\u000b4|    |\\u000b9|    |\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet i = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    const firstParagraph = paragraphs.items[0];\n    firstParagraph.insertText(newCellText[i], \"Replace\");\n    i++;\n  }\n  await context.sync();\n}\n\nawait context.sync();\n", "ps1": "# Replace the 15 lattice-multiplication exercise cells (5 rows x 3 cols)\n# with the new set of problems, keeping the table/cell/run formatting intact.\n# Each cell holds a single run with 5 \"lines\" separated by a Word line break\n# (vertical-tab char, Chr(11) / <w:br/> in the OOXML).\n$brk = [char]11\n\n$newCellLines = @(\n    @(\"60 x 73\", \"  7    3\", \"  ----\", \"6|    |\", \"0|    |\"),\n    @(\"21 x 75\", \"  7    5\", \"  ----\", \"2|    |\", \"1|    |\"),\n    @(\"58 x 47\", \"  4    7\", \"  ----\", \"5|    |\", \"8|    |\"),\n    @(\"60 x 86\", \"  8    6\", \"  ----\", \"6|    |\", \"0|    |\"),\n    @(\"56 x 42\", \"  4    2\", \"  ----\", \"5|    |\", \"6|    |\"),\n    @(\"29 x 36\", \"  3    6\", \"  ----\", \"2|    |\", \"9|    |\"),\n    @(\"20 x 89\", \"  8    9\", \"  ----\", \"2|    |\", \"0|    |\"),\n    @(\"94 x 85\", \"  8    5\", \"  ----\", \"9|    |\", \"4|    |\"),\n    @(\"88 x 40\", \"  4    0\", \"  ----\", \"8|    |\", \"8|    |\"),\n    @(\"86 x 34\", \"  3    4\", \"  ----\", \"8|    |\", \"6|    |\"),\n    @(\"15 x 55\", \"  5    5\", \"  ----\", \"1|    |\", \"5|    |\"),\n    @(\"96 x 98\", \"  9    8\", \"  ----\", \"9|    |\", \"6|    |\"),\n    @(\"43 x 24\", \"  2    4\", \"  ----\", \"4|    |\", \"3|    |\"),\n    @(\"15 x 78\", \"  7    8\", \"  ----\", \"1|    |\", \"5|    |\"),\n    @(\"49 x 42\", \"  4    2\", \"  ----\", \"4|    |\", \"9|    |\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $lines = $newCellLines[$i]\n        $text = ($lines -join $brk)\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $text\n        $i++\n    }\n}\n"}
